# Auto-generated edit script: updates recomputed market-price/profit
# columns (H, I, J, K, L, M, N) on several Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 766.087
$ws.Range("I28").Value = 744.0909
$ws.Range("K28").Value = 744.0909
$ws.Range("M28").Value = -259.0909
$ws.Range("H98").Value = 70914.734
$ws.Range("I98").Value = 37751.32
$ws.Range("K98").Value = 37751.32
$ws.Range("M98").Value = -36253.32
$ws.Range("H115").Value = 893.5333000000001
$ws.Range("I115").Value = 750.2143
$ws.Range("K115").Value = 2250.6429
$ws.Range("M115").Value = -683.6428999999998
$ws.Range("H122").Value = 70914.734
$ws.Range("I122").Value = 37751.32
$ws.Range("K122").Value = 113253.96
$ws.Range("M122").Value = -110803.96
$ws.Range("H131").Value = 2266.25
$ws.Range("I131").Value = 1688.3334
$ws.Range("K131").Value = 5065.0002
$ws.Range("M131").Value = -25.0002000000004
$ws.Range("H132").Value = 2113.3044
$ws.Range("I132").Value = 1743.1428
$ws.Range("K132").Value = 5229.428400000001
$ws.Range("M132").Value = -2699.428400000001
$ws.Range("H137").Value = 2483.4614
$ws.Range("I137").Value = 1779.6428
$ws.Range("K137").Value = 5338.928400000001
$ws.Range("M137").Value = -2788.928400000001
$ws.Range("H138").Value = 3565.0344
$ws.Range("I138").Value = 1326.2106
$ws.Range("K138").Value = 3978.6318
$ws.Range("M138").Value = 1161.3682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22065.277
$ws.Range("I32").Value = 26796.932
$ws.Range("J32").Value = 14629.821
$ws.Range("K32").Value = 26796.932
$ws.Range("L32").Value = 14629.821
$ws.Range("M32").Value = -26509.932
$ws.Range("N32").Value = -15203.821
$ws.Range("H45").Value = 6709.1816
$ws.Range("I45").Value = 7312.222
$ws.Range("K45").Value = 7312.222
$ws.Range("M45").Value = -6935.222
$ws.Range("H61").Value = 24400.377
$ws.Range("I61").Value = 2669.8518
$ws.Range("K61").Value = 2669.8518
$ws.Range("M61").Value = -2457.8518
$ws.Range("H74").Value = 64571.277
$ws.Range("I74").Value = 41300.26
$ws.Range("K74").Value = 41300.26
$ws.Range("M74").Value = -40426.26
$ws.Range("H77").Value = 64571.277
$ws.Range("I77").Value = 41300.26
$ws.Range("K77").Value = 206501.3
$ws.Range("M77").Value = -202133.3
$ws.Range("H110").Value = 29660.938
$ws.Range("I110").Value = 32362.725
$ws.Range("J110").Value = 3543.6667
$ws.Range("K110").Value = 32362.725
$ws.Range("L110").Value = 3543.6667
$ws.Range("M110").Value = -30317.725
$ws.Range("N110").Value = -7633.6667
$ws.Range("H132").Value = 2892.75
$ws.Range("I132").Value = 2671.4814
$ws.Range("K132").Value = 8014.4442
$ws.Range("M132").Value = -5484.4442
$ws.Range("H136").Value = 24400.377
$ws.Range("I136").Value = 2669.8518
$ws.Range("K136").Value = 8009.555399999999
$ws.Range("M136").Value = -5459.555399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 47500
$ws.Range("I26").Value = 47500
$ws.Range("K26").Value = 47500
$ws.Range("M26").Value = -47208
$ws.Range("H105").Value = 1737.1875
$ws.Range("I105").Value = 1504.3334
$ws.Range("J105").Value = 2994.6
$ws.Range("K105").Value = 1504.3334
$ws.Range("L105").Value = 2994.6
$ws.Range("M105").Value = 242.6666
$ws.Range("N105").Value = -6488.6
$ws.Range("H134").Value = 36848.15
$ws.Range("I134").Value = 43680.773
$ws.Range("J134").Value = 6784.6
$ws.Range("K134").Value = 131042.319
$ws.Range("L134").Value = 20353.8
$ws.Range("M134").Value = -128507.319
$ws.Range("N134").Value = -25423.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1899.9231
$ws.Range("I16").Value = 1965.2222
$ws.Range("J16").Value = 1753
$ws.Range("K16").Value = 1965.2222
$ws.Range("L16").Value = 1753
$ws.Range("M16").Value = -1678.2222
$ws.Range("N16").Value = -2327
$ws.Range("H31").Value = 2282.131
$ws.Range("I31").Value = 1957.7894
$ws.Range("J31").Value = 2818
$ws.Range("K31").Value = 1957.7894
$ws.Range("L31").Value = 2818
$ws.Range("M31").Value = -1662.7894
$ws.Range("N31").Value = -3408
$ws.Range("H34").Value = 2282.131
$ws.Range("I34").Value = 1957.7894
$ws.Range("J34").Value = 2818
$ws.Range("K34").Value = 1957.7894
$ws.Range("L34").Value = 2818
$ws.Range("M34").Value = -1755.7894
$ws.Range("N34").Value = -3222
$ws.Range("H58").Value = 9752.333000000001
$ws.Range("I58").Value = 1641.25
$ws.Range("K58").Value = 1641.25
$ws.Range("M58").Value = -1438.25
$ws.Range("H107").Value = 1006.56366
$ws.Range("I107").Value = 401.7805
$ws.Range("J107").Value = 2777.7144
$ws.Range("K107").Value = 401.7805
$ws.Range("L107").Value = 2777.7144
$ws.Range("M107").Value = 1518.2195
$ws.Range("N107").Value = -6617.7144
$ws.Range("H113").Value = 1899.9231
$ws.Range("I113").Value = 1965.2222
$ws.Range("J113").Value = 1753
$ws.Range("K113").Value = 1965.2222
$ws.Range("L113").Value = 1753
$ws.Range("M113").Value = 204.7778000000001
$ws.Range("N113").Value = -6093
$ws.Range("H132").Value = 9187.736999999999
$ws.Range("I132").Value = 4703.6
$ws.Range("J132").Value = 26003.25
$ws.Range("K132").Value = 14110.8
$ws.Range("L132").Value = 78009.75
$ws.Range("M132").Value = -11580.8
$ws.Range("N132").Value = -83069.75
$ws.Range("H134").Value = 4603.0454
$ws.Range("I134").Value = 4501.0586
$ws.Range("K134").Value = 13503.1758
$ws.Range("M134").Value = -10968.1758
$ws.Range("H136").Value = 9752.333000000001
$ws.Range("I136").Value = 1641.25
$ws.Range("K136").Value = 4923.75
$ws.Range("M136").Value = -2373.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5565.6665
$ws.Range("J42").Value = 5849
$ws.Range("L42").Value = 17547
$ws.Range("N42").Value = -18615
$ws.Range("H131").Value = 3066.6667
$ws.Range("J131").Value = 3992.2222
$ws.Range("L131").Value = 11976.6666
$ws.Range("N131").Value = -22056.6666
$ws.Range("H137").Value = 112695.375
$ws.Range("I137").Value = 1255298.6
$ws.Range("J137").Value = 5156.247
$ws.Range("K137").Value = 3765895.8
$ws.Range("L137").Value = 15468.741
$ws.Range("M137").Value = -3760795.8
$ws.Range("N137").Value = -25668.741

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 8003
$ws.Range("J19").Value = 8003
$ws.Range("L19").Value = 8003
$ws.Range("N19").Value = -8579
$ws.Range("H132").Value = 1823763.4
$ws.Range("I132").Value = 2005439.8
$ws.Range("K132").Value = 6016319.4
$ws.Range("M132").Value = -6013789.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3017.5789
$ws.Range("I61").Value = 3017.5789
$ws.Range("K61").Value = 3017.5789
$ws.Range("M61").Value = -2815.5789
$ws.Range("H113").Value = 3017.5789
$ws.Range("I113").Value = 3017.5789
$ws.Range("K113").Value = 3017.5789
$ws.Range("M113").Value = -847.5789
$ws.Range("H132").Value = 3350.225
$ws.Range("I132").Value = 3027.6177
$ws.Range("K132").Value = 9082.8531
$ws.Range("M132").Value = -6552.8531
$ws.Range("H136").Value = 3042.8
$ws.Range("I136").Value = 2429.4524
$ws.Range("K136").Value = 7288.3572
$ws.Range("M136").Value = -4738.3572
$ws.Range("H139").Value = 92825
$ws.Range("J139").Value = 105000
$ws.Range("L139").Value = 105000
$ws.Range("N139").Value = -115280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 539.9167
$ws.Range("I107").Value = 348.5
$ws.Range("K107").Value = 1045.5
$ws.Range("M107").Value = 874.5
$ws.Range("H126").Value = 26713.812
$ws.Range("I126").Value = 30244.357
$ws.Range("K126").Value = 90733.071
$ws.Range("M126").Value = -88263.071
$ws.Range("H132").Value = 5716.472
$ws.Range("I132").Value = 5212.7417
$ws.Range("J132").Value = 8839.6
$ws.Range("K132").Value = 15638.2251
$ws.Range("L132").Value = 26518.8
$ws.Range("M132").Value = -13108.2251
$ws.Range("N132").Value = -31578.8
$ws.Range("H136").Value = 4056.3333
$ws.Range("I136").Value = 4054.9443
$ws.Range("K136").Value = 12164.8329
$ws.Range("M136").Value = -9614.832900000001
